$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 20,10
$data[0,0] = -17.63924884132917
$data[0,1] = -0.6495328223752749
$data[0,2] = -17.63924884132917
$data[0,3] = -17.63924884132917
$data[0,4] = -17.63924884132917
$data[0,5] = -17.63924884132917
$data[0,6] = -17.63924884132917
$data[0,7] = -17.63924884132917
$data[0,8] = -17.63924884132917
$data[0,9] = -17.63924884132917
$data[1,0] = -17.63924884132917
$data[1,1] = -17.63924884132917
$data[1,2] = -17.63924884132917
$data[1,3] = -17.63924884132917
$data[1,4] = -17.63924884132917
$data[1,5] = -17.63924884132917
$data[1,6] = -17.63924884132917
$data[1,7] = 0.6608578978888114
$data[1,8] = -17.63924884132917
$data[1,9] = -17.63924884132917
$data[2,0] = -17.63924884132917
$data[2,1] = -0.4796829588297395
$data[2,2] = 0.1336391600099736
$data[2,3] = -17.63924884132917
$data[2,4] = 3.955418829193425
$data[2,5] = -17.63924884132917
$data[2,6] = -17.63924884132917
$data[2,7] = -17.63924884132917
$data[2,8] = 2.217469984559732
$data[2,9] = -17.63924884132917
$data[3,0] = -17.63924884132917
$data[3,1] = 0.1618342647804886
$data[3,2] = -17.63924884132917
$data[3,3] = -17.63924884132917
$data[3,4] = -17.63924884132917
$data[3,5] = 3.606935634469872
$data[3,6] = -17.63924884132917
$data[3,7] = -17.63924884132917
$data[3,8] = -17.63924884132917
$data[3,9] = -17.63924884132917
$data[4,0] = -17.63924884132917
$data[4,1] = -17.63924884132917
$data[4,2] = -17.63924884132917
$data[4,3] = -17.63924884132917
$data[4,4] = -17.63924884132917
$data[4,5] = -17.63924884132917
$data[4,6] = -17.63924884132917
$data[4,7] = -17.63924884132917
$data[4,8] = -17.63924884132917
$data[4,9] = -17.63924884132917
$data[5,0] = 3.174141723644889
$data[5,1] = -17.63924884132917
$data[5,2] = -17.63924884132917
$data[5,3] = -17.63924884132917
$data[5,4] = -17.63924884132917
$data[5,5] = -17.63924884132917
$data[5,6] = -17.63924884132917
$data[5,7] = -17.63924884132917
$data[5,8] = -17.63924884132917
$data[5,9] = -17.63924884132917
$data[6,0] = -17.63924884132917
$data[6,1] = -17.63924884132917
$data[6,2] = -17.63924884132917
$data[6,3] = 2.097927519377503
$data[6,4] = -17.63924884132917
$data[6,5] = -17.63924884132917
$data[6,6] = -17.63924884132917
$data[6,7] = -17.63924884132917
$data[6,8] = -17.63924884132917
$data[6,9] = -17.63924884132917
$data[7,0] = 3.455960790580021
$data[7,1] = -17.63924884132917
$data[7,2] = -17.63924884132917
$data[7,3] = -17.63924884132917
$data[7,4] = -17.63924884132917
$data[7,5] = -17.63924884132917
$data[7,6] = -17.63924884132917
$data[7,7] = -17.63924884132917
$data[7,8] = -17.63924884132917
$data[7,9] = -17.63924884132917
$data[8,0] = -17.63924884132917
$data[8,1] = -17.63924884132917
$data[8,2] = -17.63924884132917
$data[8,3] = -17.63924884132917
$data[8,4] = -17.63924884132917
$data[8,5] = -17.63924884132917
$data[8,6] = -17.63924884132917
$data[8,7] = 0.4897372735290023
$data[8,8] = -17.63924884132917
$data[8,9] = 1.413205114717469
$data[9,0] = -17.63924884132917
$data[9,1] = -17.63924884132917
$data[9,2] = -17.63924884132917
$data[9,3] = 1.916637771081118
$data[9,4] = -17.63924884132917
$data[9,5] = 1.695829606417858
$data[9,6] = -17.63924884132917
$data[9,7] = -17.63924884132917
$data[9,8] = -17.63924884132917
$data[9,9] = 1.195793055499914
$data[10,0] = -17.63924884132917
$data[10,1] = -17.63924884132917
$data[10,2] = -17.63924884132917
$data[10,3] = -17.63924884132917
$data[10,4] = -17.63924884132917
$data[10,5] = -17.63924884132917
$data[10,6] = -17.63924884132917
$data[10,7] = -17.63924884132917
$data[10,8] = -17.63924884132917
$data[10,9] = -17.63924884132917
$data[11,0] = -17.63924884132917
$data[11,1] = -17.63924884132917
$data[11,2] = -17.63924884132917
$data[11,3] = 1.819752443810723
$data[11,4] = -17.63924884132917
$data[11,5] = -17.63924884132917
$data[11,6] = -17.63924884132917
$data[11,7] = -17.63924884132917
$data[11,8] = 0.5902318717309779
$data[11,9] = 2.157812597682025
$data[12,0] = -17.63924884132917
$data[12,1] = -17.63924884132917
$data[12,2] = 1.825865135848753
$data[12,3] = -17.63924884132917
$data[12,4] = -17.63924884132917
$data[12,5] = -17.63924884132917
$data[12,6] = -17.63924884132917
$data[12,7] = -17.63924884132917
$data[12,8] = -17.63924884132917
$data[12,9] = 1.801014075638085
$data[13,0] = -17.63924884132917
$data[13,1] = -17.63924884132917
$data[13,2] = -0.5785845702558944
$data[13,3] = -17.63924884132917
$data[13,4] = -17.63924884132917
$data[13,5] = -17.63924884132917
$data[13,6] = -17.63924884132917
$data[13,7] = -17.63924884132917
$data[13,8] = -17.63924884132917
$data[13,9] = -17.63924884132917
$data[14,0] = -17.63924884132917
$data[14,1] = -17.63924884132917
$data[14,2] = -17.63924884132917
$data[14,3] = -17.63924884132917
$data[14,4] = -17.63924884132917
$data[14,5] = -17.63924884132917
$data[14,6] = -17.63924884132917
$data[14,7] = -17.63924884132917
$data[14,8] = 3.022617091466146
$data[14,9] = -17.63924884132917
$data[15,0] = -17.63924884132917
$data[15,1] = 0.02531649626609556
$data[15,2] = -0.4392602515516793
$data[15,3] = -17.63924884132917
$data[15,4] = -17.63924884132917
$data[15,5] = -17.63924884132917
$data[15,6] = 4.321921381276691
$data[15,7] = -0.8816374648506528
$data[15,8] = 1.241939819526277
$data[15,9] = -17.63924884132917
$data[16,0] = -17.63924884132917
$data[16,1] = -17.63924884132917
$data[16,2] = -17.63924884132917
$data[16,3] = -17.63924884132917
$data[16,4] = -17.63924884132917
$data[16,5] = -17.63924884132917
$data[16,6] = -17.63924884132917
$data[16,7] = -1.05286348475474
$data[16,8] = 1.745064853649958
$data[16,9] = -17.63924884132917
$data[17,0] = -17.63924884132917
$data[17,1] = -17.63924884132917
$data[17,2] = 2.921290178741247
$data[17,3] = -17.63924884132917
$data[17,4] = -17.63924884132917
$data[17,5] = -17.63924884132917
$data[17,6] = -17.63924884132917
$data[17,7] = 2.421564492404332
$data[17,8] = -17.63924884132917
$data[17,9] = -17.63924884132917
$data[18,0] = -17.63924884132917
$data[18,1] = 3.069620340633934
$data[18,2] = 2.672496148222993
$data[18,3] = -17.63924884132917
$data[18,4] = 2.165677805333337
$data[18,5] = -17.63924884132917
$data[18,6] = -17.63924884132917
$data[18,7] = 3.41040098892621
$data[18,8] = -17.63924884132917
$data[18,9] = 2.827597344485611
$data[19,0] = -17.63924884132917
$data[19,1] = 3.020302026574777
$data[19,2] = -17.63924884132917
$data[19,3] = 3.07268106595233
$data[19,4] = -17.63924884132917
$data[19,5] = 2.194124622555601
$data[19,6] = -17.63924884132917
$data[19,7] = -17.63924884132917
$data[19,8] = -17.63924884132917
$data[19,9] = -17.63924884132917

$ws.Range("B2:K21").Value = $data
